$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.202.46"
$ws.Range("E2").Value = "  +5.11%  "
Set-TextValue $ws.Range("D3") "3.457.83"
$ws.Range("E3").Value = "  +4.64%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "582.24"
$ws.Range("E5").Value = "  +5.86%  "
Set-TextValue $ws.Range("D6") "184.22"
$ws.Range("E6").Value = "  +7.16%  "
$ws.Range("E7").Value = "  +2.66%  "
Set-TextValue $ws.Range("D8") "3.451.38"
$ws.Range("E8").Value = "  +4.61%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +2.33%  "
Set-TextValue $ws.Range("D11") "0.644"
$ws.Range("E11").Value = "  +2.85%  "
Set-TextValue $ws.Range("D12") "56.27"
$ws.Range("E12").Value = "  +5.79%  "
Set-TextValue $ws.Range("D13") "0.0000278"
$ws.Range("E13").Value = "  +0.15%  "
Set-TextValue $ws.Range("D14") "9.43"
$ws.Range("E14").Value = "  +5.02%  "
Set-TextValue $ws.Range("D15") "4.004.13"
$ws.Range("E15").Value = "  +4.40%  "
Set-TextValue $ws.Range("D16") "18.63"
$ws.Range("E16").Value = "  +3.51%  "
Set-TextValue $ws.Range("D17") "3.450.78"
$ws.Range("E17").Value = "  +4.26%  "
Set-TextValue $ws.Range("D18") "67.077.89"
$ws.Range("E18").Value = "  +4.91%  "
Set-TextValue $ws.Range("D20") "12.09"
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("E21").Value = "  +3.79%  "
Set-TextValue $ws.Range("D22") "482.23"
$ws.Range("E22").Value = "  +5.30%  "
Set-TextValue $ws.Range("D23") "5.47"
$ws.Range("E23").Value = "  +9.72%  "
Set-TextValue $ws.Range("D24") "16.74"
$ws.Range("E24").Value = "  +21.49%  "
Set-TextValue $ws.Range("D25") "4.43"
$ws.Range("E25").Value = "  +9.84%  "
Set-TextValue $ws.Range("D26") "89.66"
$ws.Range("E26").Value = "  +3.73%  "
Set-TextValue $ws.Range("D27") "2.94"
$ws.Range("E27").Value = "  +3.50%  "
Set-TextValue $ws.Range("D28") "10.97"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("E29").Value = "  +6.83%  "
Set-TextValue $ws.Range("D30") "31.34"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("E31").Value = "  +9.60%  "
Set-TextValue $ws.Range("D32") "64.25"
$ws.Range("E32").Value = "  +5.38%  "
Set-TextValue $ws.Range("D33") "11.71"
$ws.Range("E33").Value = "  +3.28%  "
Set-TextValue $ws.Range("D34") "588.33"
$ws.Range("E34").Value = "  +4.80%  "
Set-TextValue $ws.Range("D35") "0.112"
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("E37").Value = "  -0.01%  "
Set-TextValue $ws.Range("D38") "3.56"
$ws.Range("E38").Value = "  +2.04%  "
Set-TextValue $ws.Range("D39") "36.44"
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$pepePrice = [string]::Concat("0.0", [char]0x2083, "0769")
Set-TextValue $ws.Range("D40") $pepePrice
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D41") "0.384"
$ws.Range("E41").Value = "  +5.81%  "
Set-TextValue $ws.Range("D42") "3.211.09"
$ws.Range("E42").Value = "  +5.89%  "
Set-TextValue $ws.Range("D43") "2.90"
$ws.Range("E43").Value = "  +6.19%  "
$ws.Range("E44").Value = "  +4.46%  "
Set-TextValue $ws.Range("D45") "2.53"
$ws.Range("E45").Value = "  +4.51%  "
$ws.Range("E46").Value = "  +21.95%  "
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E48").Value = "  +2.19%  "
Set-TextValue $ws.Range("D49") "8.75"
$ws.Range("E49").Value = "  +8.37%  "
Set-TextValue $ws.Range("D50") "0.999"
Set-TextValue $ws.Range("D51") "3.21"
$ws.Range("E51").Value = "  +11.58%  "

Write-Output "done"
